$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), reusing the same formatting as
# the existing header row (bold font + border), matching style of H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF)
$data = @(
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(8, 9),
    @(8, 8),
    @(12, 12),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(6, 6),
    @(5, 5),
    @(5, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
